# Auto-generated edit script applying numeric corrections to the
# per-profession Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR),
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4773.75
$ws.Range("I2").Value = 1658.2
$ws.Range("K2").Value = 1658.2
$ws.Range("M2").Value = -1545.2
$ws.Range("H76").Value = 3177
$ws.Range("I76").Value = 3065.1765
$ws.Range("J76").Value = 4444.3335
$ws.Range("K76").Value = 3065.1765
$ws.Range("L76").Value = 4444.3335
$ws.Range("M76").Value = -2750.1765
$ws.Range("N76").Value = -5074.3335
$ws.Range("H79").Value = 3177
$ws.Range("I79").Value = 3065.1765
$ws.Range("J79").Value = 4444.3335
$ws.Range("K79").Value = 3065.1765
$ws.Range("L79").Value = 4444.3335
$ws.Range("M79").Value = -1973.1765
$ws.Range("N79").Value = -6628.3335
$ws.Range("H111").Value = 1146.5
$ws.Range("J111").Value = 1159.3334
$ws.Range("L111").Value = 3478.0002
$ws.Range("N111").Value = -9612.0002
$ws.Range("H112").Value = 2420.5454
$ws.Range("I112").Value = 2248.8333
$ws.Range("J112").Value = 2458.7036
$ws.Range("K112").Value = 6746.499899999999
$ws.Range("L112").Value = 7376.110799999999
$ws.Range("M112").Value = -5638.499899999999
$ws.Range("N112").Value = -9592.110799999999
$ws.Range("H125").Value = 1328.6364
$ws.Range("I125").Value = 1047.875
$ws.Range("K125").Value = 9430.875
$ws.Range("M125").Value = -6970.875
$ws.Range("H133").Value = 199998.4
$ws.Range("H135").Value = 590.94116
$ws.Range("I135").Value = 552.875
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 4975.875
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -2440.875
$ws.Range("N135").Value = -15870
$ws.Range("H138").Value = 4901.093
$ws.Range("I138").Value = 5450.24
$ws.Range("J138").Value = 4427.6895
$ws.Range("K138").Value = 16350.72
$ws.Range("L138").Value = 13283.0685
$ws.Range("M138").Value = -11210.72
$ws.Range("N138").Value = -23563.0685
$ws.Range("H141").Value = 9750
$ws.Range("I141").Value = 4500
$ws.Range("J141").Value = 11500
$ws.Range("K141").Value = 13500
$ws.Range("L141").Value = 34500
$ws.Range("M141").Value = -8320
$ws.Range("N141").Value = -44860

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1053.5714
$ws.Range("I5").Value = 895.8333
$ws.Range("K5").Value = 895.8333
$ws.Range("M5").Value = -783.8333
$ws.Range("H32").Value = 3052.3635
$ws.Range("I32").Value = 1907.2264
$ws.Range("K32").Value = 1907.2264
$ws.Range("M32").Value = -1620.2264
$ws.Range("H45").Value = 80517.84
$ws.Range("I45").Value = 80517.84
$ws.Range("K45").Value = 80517.84
$ws.Range("M45").Value = -80140.84
$ws.Range("H61").Value = 1764590.9
$ws.Range("J61").Value = 3681089
$ws.Range("L61").Value = 3681089
$ws.Range("N61").Value = -3681513
$ws.Range("H92").Value = 42498.5
$ws.Range("J92").Value = 42498.5
$ws.Range("L92").Value = 42498.5
$ws.Range("N92").Value = -47490.5
$ws.Range("H102").Value = 2980
$ws.Range("I102").Value = 2980
$ws.Range("K102").Value = 2980
$ws.Range("M102").Value = -1358
$ws.Range("H136").Value = 1764590.9
$ws.Range("J136").Value = 3681089
$ws.Range("L136").Value = 11043267
$ws.Range("N136").Value = -11048367

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1053.5714
$ws.Range("I4").Value = 895.8333
$ws.Range("K4").Value = 895.8333
$ws.Range("M4").Value = -780.8333
$ws.Range("H99").Value = 13276.9375
$ws.Range("I99").Value = 12639.333
$ws.Range("K99").Value = 12639.333
$ws.Range("M99").Value = -11141.333

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 18209.7
$ws.Range("I86").Value = 6774.25
$ws.Range("K86").Value = 6774.25
$ws.Range("M86").Value = -5651.25
$ws.Range("H89").Value = 18209.7
$ws.Range("I89").Value = 6774.25
$ws.Range("K89").Value = 33871.25
$ws.Range("M89").Value = -28255.25
$ws.Range("H134").Value = 3100.1892
$ws.Range("I134").Value = 2690.7585
$ws.Range("K134").Value = 8072.2755
$ws.Range("M134").Value = -5537.2755

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 661.3461
$ws.Range("I2").Value = 330.875
$ws.Range("K2").Value = 1985.25
$ws.Range("M2").Value = -1872.25
$ws.Range("H68").Value = 1299.4054
$ws.Range("J68").Value = 1223.7587
$ws.Range("L68").Value = 3671.2761
$ws.Range("N68").Value = -5293.2761
$ws.Range("H71").Value = 1299.4054
$ws.Range("J71").Value = 1223.7587
$ws.Range("L71").Value = 11013.8283
$ws.Range("N71").Value = -19125.8283
$ws.Range("H82").Value = 13073.833
$ws.Range("I82").Value = 3500
$ws.Range("J82").Value = 14988.6
$ws.Range("K82").Value = 10500
$ws.Range("L82").Value = 44965.8
$ws.Range("M82").Value = -10094
$ws.Range("N82").Value = -45777.8
$ws.Range("H85").Value = 13073.833
$ws.Range("I85").Value = 3500
$ws.Range("J85").Value = 14988.6
$ws.Range("K85").Value = 10500
$ws.Range("L85").Value = 44965.8
$ws.Range("M85").Value = -9096
$ws.Range("N85").Value = -47773.8
$ws.Range("H120").Value = 23749
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 23749
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 71247
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -80923
$ws.Range("H139").Value = 5880.8335
$ws.Range("I139").Value = 4958.2
$ws.Range("K139").Value = 14874.6
$ws.Range("M139").Value = -9734.599999999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 572071.4399999999
$ws.Range("J3").Value = 667400
$ws.Range("L3").Value = 667400
$ws.Range("N3").Value = -667632
$ws.Range("H13").Value = 2321.1667
$ws.Range("H52").Value = 4999.875
$ws.Range("I52").Value = 4999
$ws.Range("K52").Value = 4999
$ws.Range("M52").Value = -4740
$ws.Range("H69").Value = 60000
$ws.Range("J69").Value = 60000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61498
$ws.Range("H72").Value = 60000
$ws.Range("J72").Value = 60000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -187488
$ws.Range("H80").Value = 147359.7
$ws.Range("I80").Value = 151117.16
$ws.Range("J80").Value = 137339.83
$ws.Range("K80").Value = 151117.16
$ws.Range("L80").Value = 137339.83
$ws.Range("M80").Value = -150119.16
$ws.Range("N80").Value = -139335.83
$ws.Range("H83").Value = 147359.7
$ws.Range("I83").Value = 151117.16
$ws.Range("J83").Value = 137339.83
$ws.Range("K83").Value = 755585.8
$ws.Range("L83").Value = 686699.1499999999
$ws.Range("M83").Value = -750593.8
$ws.Range("N83").Value = -696683.1499999999
$ws.Range("H132").Value = 11029681
$ws.Range("I132").Value = 3968.4443
$ws.Range("K132").Value = 11905.3329
$ws.Range("M132").Value = -9375.332900000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2884.4092
$ws.Range("I16").Value = 2688.4285
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 2688.4285
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = -2518.4285
$ws.Range("N16").Value = -7340
$ws.Range("H20").Value = 110106110
$ws.Range("J20").Value = 157150160
$ws.Range("L20").Value = 157150160
$ws.Range("N20").Value = -157150612
$ws.Range("H35").Value = 1183.3334
$ws.Range("I35").Value = 1183.3334
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1183.3334
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -847.3334
$ws.Range("N35").ClearContents()
$ws.Range("H100").Value = 4689.8887
$ws.Range("I100").Value = 3852.4
$ws.Range("K100").Value = 3852.4
$ws.Range("M100").Value = -3311.4

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 3500774.8
$ws.Range("J5").Value = 1334366.4
$ws.Range("L5").Value = 1334366.4
$ws.Range("N5").Value = -1334590.4
$ws.Range("H23").Value = 10494.5
$ws.Range("J23").Value = 19989
$ws.Range("L23").Value = 19989
$ws.Range("N23").Value = -20447
$ws.Range("H107").Value = 1593.4642
$ws.Range("I107").Value = 1185.05
$ws.Range("J107").Value = 2614.5
$ws.Range("K107").Value = 3555.15
$ws.Range("L107").Value = 7843.5
$ws.Range("M107").Value = -1635.15
$ws.Range("N107").Value = -11683.5
$ws.Range("H132").Value = 2113.5
$ws.Range("I132").Value = 1837.3529
$ws.Range("J132").Value = 2784.1428
$ws.Range("K132").Value = 5512.0587
$ws.Range("L132").Value = 8352.428400000001
$ws.Range("M132").Value = -2982.0587
$ws.Range("N132").Value = -13412.4284
